# Auto-generated script applying cell updates described in the commit
# 'Add data for 2024-08-30' - updates 2024 YTD (column K, and a couple
# column I 2022 corrections) violent-crime counts across the citywide
# totals sheet, the by-neighborhood rollup sheet, and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('K2').Value = 5197
$ws.Range('K3').Value = 5340
$ws.Range('I4').Value = 1224
$ws.Range('K4').Value = 1115
$ws.Range('K5').Value = 383
$ws.Range('K6').Value = 5972
$ws.Range('I7').Value = 16842
$ws.Range('K7').Value = 18007

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range('K4').Value = 8
$ws.Range('K7').Value = 40

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('K2').Value = 332
$ws.Range('K3').Value = 366
$ws.Range('K6').Value = 414
$ws.Range('K7').Value = 1219

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('K3').Value = 145
$ws.Range('K7').Value = 402

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('K2').Value = 208
$ws.Range('K3').Value = 284
$ws.Range('K6').Value = 221
$ws.Range('K7').Value = 765

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('K2').Value = 103
$ws.Range('K3').Value = 111
$ws.Range('K6').Value = 71
$ws.Range('K7').Value = 310

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('K2').Value = 175
$ws.Range('K3').Value = 202
$ws.Range('K7').Value = 609

$ws = $wb.Worksheets.Item('New City')
$ws.Range('K2').Value = 128
$ws.Range('K3').Value = 102
$ws.Range('K7').Value = 409

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('K3').Value = 126
$ws.Range('K7').Value = 306

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('K2').Value = 156
$ws.Range('K7').Value = 536
$ws.Range('K8').Value = 1219
$ws.Range('K11').Value = 346
$ws.Range('K13').Value = 21
$ws.Range('K16').Value = 55
$ws.Range('K20').Value = 415
$ws.Range('K23').Value = 185
$ws.Range('K27').Value = 168
$ws.Range('K29').Value = 963
$ws.Range('K31').Value = 194
$ws.Range('K33').Value = 765
$ws.Range('K34').Value = 101
$ws.Range('K36').Value = 236
$ws.Range('K37').Value = 609
$ws.Range('K42').Value = 660
$ws.Range('K48').Value = 224
$ws.Range('K52').Value = 472
$ws.Range('K54').Value = 354
$ws.Range('K55').Value = 201
$ws.Range('K57').Value = 66
$ws.Range('I63').Value = 145
$ws.Range('K64').Value = 113
$ws.Range('K65').Value = 409
$ws.Range('K67').Value = 685
$ws.Range('K68').Value = 48
$ws.Range('K69').Value = 40
$ws.Range('K71').Value = 57
$ws.Range('K72').Value = 87
$ws.Range('K73').Value = 154
$ws.Range('K78').Value = 204
$ws.Range('K79').Value = 445
$ws.Range('K81').Value = 13
$ws.Range('K83').Value = 402
$ws.Range('K85').Value = 846
$ws.Range('K86').Value = 121
$ws.Range('K89').Value = 264
$ws.Range('K91').Value = 194
$ws.Range('K92').Value = 68
$ws.Range('K95').Value = 310
$ws.Range('K96').Value = 193
$ws.Range('K97').Value = 143
$ws.Range('K98').Value = 85
$ws.Range('K99').Value = 306
$ws.Range('K100').Value = 34
$ws.Range('I101').Value = 16842
$ws.Range('K101').Value = 18007

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('K2').Value = 67
$ws.Range('K7').Value = 194

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('K3').Value = 243
$ws.Range('K6').Value = 191
$ws.Range('K7').Value = 685

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('K2').Value = 56
$ws.Range('K6').Value = 188
$ws.Range('K7').Value = 354

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('K3').Value = 346
$ws.Range('K7').Value = 963

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('K4').Value = 33
$ws.Range('K7').Value = 224

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('K3').Value = 204
$ws.Range('K6').Value = 251
$ws.Range('K7').Value = 660

$ws = $wb.Worksheets.Item('Boystown')
$ws.Range('K4').Value = 4
$ws.Range('K6').Value = 21

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('K2').Value = 60
$ws.Range('K7').Value = 204

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('K2').Value = 61
$ws.Range('K7').Value = 201

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('K3').Value = 67
$ws.Range('K7').Value = 185

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('K4').Value = 10
$ws.Range('K7').Value = 193

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('K3').Value = 92
$ws.Range('K6').Value = 46
$ws.Range('K7').Value = 194

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('K3').Value = 144
$ws.Range('K7').Value = 445

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('K4').Value = 14
$ws.Range('K7').Value = 113

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('K3').Value = 133
$ws.Range('K5').Value = 7
$ws.Range('K7').Value = 415

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('K2').Value = 93
$ws.Range('K3').Value = 67
$ws.Range('K7').Value = 236

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range('K2').Value = 7
$ws.Range('K7').Value = 34

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('K6').Value = 141
$ws.Range('K7').Value = 536

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('K3').Value = 26
$ws.Range('K7').Value = 101

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range('K4').Value = 4
$ws.Range('K7').Value = 85

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('K2').Value = 116
$ws.Range('K7').Value = 346

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('K2').Value = 48
$ws.Range('K7').Value = 154

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('K2').Value = 52
$ws.Range('K7').Value = 156

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('K3').Value = 26
$ws.Range('K7').Value = 143

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range('K6').Value = 34
$ws.Range('K7').Value = 68

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('K2').Value = 74
$ws.Range('K7').Value = 264

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('K3').Value = 39
$ws.Range('K7').Value = 168

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('K4').Value = 52
$ws.Range('K7').Value = 121

$ws = $wb.Worksheets.Item('North Park')
$ws.Range('K3').Value = 10
$ws.Range('K7').Value = 48

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range('K2').Value = 18
$ws.Range('K7').Value = 66

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('K2').Value = 284
$ws.Range('K3').Value = 284
$ws.Range('K6').Value = 202
$ws.Range('K7').Value = 846

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range('K2').Value = 22
$ws.Range('K7').Value = 57

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range('K2').Value = 14
$ws.Range('K6').Value = 45
$ws.Range('K7').Value = 87

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('K2').Value = 126
$ws.Range('K3').Value = 131
$ws.Range('K7').Value = 472

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range('K6').Value = 32
$ws.Range('K7').Value = 55

$ws = $wb.Worksheets.Item('Sauganash,Forest Glen')
$ws.Range('K2').Value = 5
$ws.Range('K7').Value = 13
